$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Velutinae" subsection to "Coccineae" for the affected species
# rows (buckleyi, acerifolia, shumardii, velutina, coccinea, ellipsoidalis,
# rubra) in both the "subsection" (D) and "map" (E) columns.
$ws.Range("D26:D32").Value = "Coccineae"
$ws.Range("E26:E32").Value = "Coccineae"

# Turn on the AutoFilter for the map column (E1:E183) and keep the
# worksheet's hidden _FilterDatabase name in sync with it.
$ws.Range("E1:E183").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$E`$1:`$E`$183"
    }
}

# Update the view: scroll back to the top and select I27.
$ws.Activate()
$ws.Range("I27").Select()
